$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, pushing existing rows 66-72 down to 67-73
$ws.Rows.Item(66).Insert()

# Populate the new row 66 with the latest weekly entry
$ws.Cells.Item(66, 1).Value = 2
$ws.Cells.Item(66, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = 44826
$ws.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(66, 5).Value = 4
$ws.Cells.Item(66, 6).Value = 100112026
$ws.Cells.Item(66, 7).Value = "Haba"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 1200
$ws.Cells.Item(66, 11).Value = 7000
$ws.Cells.Item(66, 12).Value = 8000
$ws.Cells.Item(66, 13).Value = 7500
$ws.Cells.Item(66, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(66, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(66, 16).Value = 300
$ws.Cells.Item(66, 17).Value = 25
$ws.Cells.Item(66, 18).Value = "Hortaliza"
